$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4 (shifts old row4.. down by one)
$ws.Rows.Item(4).Insert()

# Copy formatting from the single cell A5 (Inzynierka row - matches target height/style) into A4 only
$ws.Range("A5").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Rows.Item(4).RowHeight = 30.75

Write-Output "done"
